# Update "想去人数" (want-to-go count) figures by +1 on two sheets,
# matching the upstream data refresh recorded in the commit.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 247
$ws1.Range("F4").Value = 871

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 247
$ws4.Range("F5").Value = 871
